# Update gh-pages to output generated at 456a3b4
# This updates the "want to go" counts (column F) and "lowest price" (column G)
# across the 展览 (Exhibition), 演出 (Performance) and 全部类型 (All types) sheets.

$wb = $excel.ActiveWorkbook

# ---- Sheet: 展览 (Exhibition) ----
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 2039
$ws1.Range("F3").Value = 659
$ws1.Range("G3").Value = 65
$ws1.Range("F4").Value = 1305
$ws1.Range("F6").Value = 68
$ws1.Range("F7").Value = 155
$ws1.Range("F11").Value = 943
$ws1.Range("F12").Value = 299
$ws1.Range("F13").Value = 160
$ws1.Range("F14").Value = 40
$ws1.Range("F17").Value = 319
$ws1.Range("F18").Value = 735
$ws1.Range("F19").Value = 114
$ws1.Range("F20").Value = 694
$ws1.Range("F21").Value = 236
$ws1.Range("F22").Value = 60
$ws1.Range("F23").Value = 952
$ws1.Range("F24").Value = 412
$ws1.Range("F25").Value = 228
$ws1.Range("F26").Value = 74
$ws1.Range("F27").Value = 338
$ws1.Range("F29").Value = 30

# ---- Sheet: 演出 (Performance) ----
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F4").Value = 347
$ws2.Range("F6").Value = 37

# ---- Sheet: 全部类型 (All types) ----
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 2039
$ws4.Range("F4").Value = 659
$ws4.Range("G4").Value = 65
$ws4.Range("F5").Value = 1305
$ws4.Range("F8").Value = 68
$ws4.Range("F9").Value = 155
$ws4.Range("F13").Value = 943
$ws4.Range("F14").Value = 299
$ws4.Range("F15").Value = 160
$ws4.Range("F17").Value = 40
$ws4.Range("F18").Value = 347
$ws4.Range("F21").Value = 354
$ws4.Range("F22").Value = 37
$ws4.Range("F24").Value = 319
$ws4.Range("F25").Value = 735
$ws4.Range("F26").Value = 114
$ws4.Range("F27").Value = 694
$ws4.Range("F28").Value = 236
$ws4.Range("F29").Value = 60
$ws4.Range("F30").Value = 952
$ws4.Range("F31").Value = 412
$ws4.Range("F34").Value = 228
$ws4.Range("F35").Value = 74
$ws4.Range("F36").Value = 338
$ws4.Range("F40").Value = 30
